$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 38501.383
$ws.Range("I11").Value = 38501.383
$ws.Range("K11").Value = 38501.383
$ws.Range("M11").Value = -38361.383
$ws.Range("H18").Value = 499
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 498
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 498
$ws.Range("M18").Value = -216
$ws.Range("N18").Value = -1066
$ws.Range("H37").Value = 938
$ws.Range("I37").Value = 957
$ws.Range("K37").Value = 2871
$ws.Range("M37").Value = -2745
$ws.Range("H45").Value = 5016.5
$ws.Range("J45").Value = 5016
$ws.Range("L45").Value = 15048
$ws.Range("N45").Value = -15432
$ws.Range("H49").Value = 2702.3333
$ws.Range("I49").Value = 1239
$ws.Range("K49").Value = 3717
$ws.Range("M49").Value = -3581
$ws.Range("H51").Value = 24748.5
$ws.Range("I51").Value = 85000
$ws.Range("J51").Value = 4664.6665
$ws.Range("K51").Value = 85000
$ws.Range("L51").Value = 4664.6665
$ws.Range("M51").Value = -84516
$ws.Range("N51").Value = -5632.6665
$ws.Range("H69").Value = 3663
$ws.Range("H72").Value = 3663
$ws.Range("H80").Value = 5382.385
$ws.Range("J80").Value = 5397
$ws.Range("L80").Value = 16191
$ws.Range("N80").Value = -18187
$ws.Range("H83").Value = 5382.385
$ws.Range("J83").Value = 5397
$ws.Range("L83").Value = 48573
$ws.Range("N83").Value = -58557
$ws.Range("H86").Value = 10866.546
$ws.Range("I86").Value = 3266.1667
$ws.Range("K86").Value = 3266.1667
$ws.Range("M86").Value = -2143.1667
$ws.Range("H89").Value = 10866.546
$ws.Range("I89").Value = 3266.1667
$ws.Range("K89").Value = 16330.8335
$ws.Range("M89").Value = -10714.8335
$ws.Range("H112").Value = 1148.125
$ws.Range("J112").Value = 1265.3334
$ws.Range("L112").Value = 3796.0002
$ws.Range("N112").Value = -6012.0002
$ws.Range("H116").Value = 45602.04
$ws.Range("I116").Value = 80581.30499999999
$ws.Range("K116").Value = 80581.30499999999
$ws.Range("M116").Value = -77139.30499999999
$ws.Range("H127").Value = 1917.2
$ws.Range("I127").Value = 2024
$ws.Range("J127").Value = 1490
$ws.Range("K127").Value = 6072
$ws.Range("L127").Value = 4470
$ws.Range("M127").Value = -1112
$ws.Range("N127").Value = -14390
$ws.Range("H132").Value = 48407.316
$ws.Range("I132").Value = 53008.85
$ws.Range("K132").Value = 159026.55
$ws.Range("M132").Value = -156496.55

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1795.9286
$ws.Range("I45").Value = 1762
$ws.Range("J45").Value = 1999.5
$ws.Range("K45").Value = 1762
$ws.Range("L45").Value = 1999.5
$ws.Range("M45").Value = -1385
$ws.Range("N45").Value = -2753.5
$ws.Range("H132").Value = 16672466
$ws.Range("I132").Value = 4205
$ws.Range("K132").Value = 12615
$ws.Range("M132").Value = -10085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1189.3043
$ws.Range("I20").Value = 1078.5555
$ws.Range("K20").Value = 1078.5555
$ws.Range("M20").Value = -831.5554999999999
$ws.Range("H86").Value = 1632.6
$ws.Range("I86").Value = 1632.6
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1632.6
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -509.5999999999999
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1632.6
$ws.Range("I89").Value = 1632.6
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 8163
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -2547
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 29549.834
$ws.Range("J107").Value = 64394.75
$ws.Range("L107").Value = 64394.75
$ws.Range("N107").Value = -68234.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 540.35297
$ws.Range("J94").Value = 462.3125
$ws.Range("L94").Value = 462.3125
$ws.Range("N94").Value = -1364.3125
$ws.Range("H107").Value = 6058.8
$ws.Range("I107").Value = 4259.3335
$ws.Range("K107").Value = 4259.3335
$ws.Range("M107").Value = -2339.3335
$ws.Range("H134").Value = 2480.5
$ws.Range("I134").Value = 2478.7273
$ws.Range("K134").Value = 7436.1819
$ws.Range("M134").Value = -4901.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H70").Value = 4887.364
$ws.Range("I70").Value = 4782.5557
$ws.Range("K70").Value = 4782.5557
$ws.Range("M70").Value = -4512.5557
$ws.Range("H73").Value = 4887.364
$ws.Range("I73").Value = 4782.5557
$ws.Range("K73").Value = 4782.5557
$ws.Range("M73").Value = -3846.5557
$ws.Range("H80").Value = 2585.2222
$ws.Range("I80").Value = 2033.5
$ws.Range("J80").Value = 6999
$ws.Range("K80").Value = 2033.5
$ws.Range("L80").Value = 6999
$ws.Range("M80").Value = -1035.5
$ws.Range("N80").Value = -8995
$ws.Range("H83").Value = 2585.2222
$ws.Range("I83").Value = 2033.5
$ws.Range("J83").Value = 6999
$ws.Range("K83").Value = 10167.5
$ws.Range("L83").Value = 34995
$ws.Range("M83").Value = -5175.5
$ws.Range("N83").Value = -44979
$ws.Range("H102").Value = 2569.1924
$ws.Range("I102").Value = 1980.0667
$ws.Range("K102").Value = 1980.0667
$ws.Range("M102").Value = -358.0667000000001
$ws.Range("H122").Value = 2813.2632
$ws.Range("I122").Value = 3203
$ws.Range("J122").Value = 1722
$ws.Range("K122").Value = 9609
$ws.Range("L122").Value = 5166
$ws.Range("M122").Value = -7159
$ws.Range("N122").Value = -10066
$ws.Range("H126").Value = 2440.7646
$ws.Range("I126").Value = 1870.4286
$ws.Range("K126").Value = 5611.2858
$ws.Range("M126").Value = -3141.2858
$ws.Range("H132").Value = 2898.1765
$ws.Range("I132").Value = 2439
$ws.Range("J132").Value = 3740
$ws.Range("K132").Value = 7317
$ws.Range("L132").Value = 11220
$ws.Range("M132").Value = -4787
$ws.Range("N132").Value = -16280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1617.5
$ws.Range("J16").Value = 1396
$ws.Range("L16").Value = 1396
$ws.Range("N16").Value = -1736
$ws.Range("H40").Value = 3915.45
$ws.Range("I40").Value = 2836.5
$ws.Range("K40").Value = 2836.5
$ws.Range("M40").Value = -2700.5
$ws.Range("H46").Value = 1295.6428
$ws.Range("I46").Value = 439.4
$ws.Range("J46").Value = 1481.7826
$ws.Range("K46").Value = 439.4
$ws.Range("L46").Value = 1481.7826
$ws.Range("M46").Value = -251.4
$ws.Range("N46").Value = -1857.7826
$ws.Range("H61").Value = 2900
$ws.Range("I61").Value = 2888.875
$ws.Range("K61").Value = 2888.875
$ws.Range("M61").Value = -2686.875
$ws.Range("H113").Value = 2900
$ws.Range("I113").Value = 2888.875
$ws.Range("K113").Value = 2888.875
$ws.Range("M113").Value = -718.875
$ws.Range("H132").Value = 3903.75
$ws.Range("I132").Value = 3061
$ws.Range("K132").Value = 9183
$ws.Range("M132").Value = -6653

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 37099
$ws.Range("J81").Value = 36549.5
$ws.Range("L81").Value = 73099
$ws.Range("N81").Value = -75221
$ws.Range("H84").Value = 37099
$ws.Range("J84").Value = 36549.5
$ws.Range("L84").Value = 365495
$ws.Range("N84").Value = -376103
$ws.Range("H113").Value = 3767.9167
$ws.Range("I113").Value = 767.4
$ws.Range("J113").Value = 5911.143
$ws.Range("K113").Value = 2302.2
$ws.Range("L113").Value = 17733.429
$ws.Range("M113").Value = -132.1999999999998
$ws.Range("N113").Value = -22073.429
$ws.Range("H124").Value = 38999.5
$ws.Range("J124").Value = 38999.5
$ws.Range("L124").Value = 38999.5
$ws.Range("N124").Value = -48819.5
$ws.Range("H132").Value = 781.8
$ws.Range("I132").Value = 757.55554
$ws.Range("K132").Value = 2272.66662
$ws.Range("M132").Value = 257.33338
$ws.Range("H136").Value = 615
$ws.Range("I136").Value = 615
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 1845
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 705
$ws.Range("N136").ClearContents()
